# Updates cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.986.42"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").Value = "2.473.77"
$ws.Range("E3").Value = "  -1.71%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.46"
$ws.Range("E5").Value = "  -1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.08"
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.508"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.156"
$ws.Range("E9").Value = "  -3.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  -1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.345"
$ws.Range("E11").Value = "  -4.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.83"
$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").Value = "2.928.70"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "68.892.57"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.90"
$ws.Range("E16").Value = "  -4.06%  "

$ws.Range("D17").Value = "2.502.41"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.88"
$ws.Range("E18").Value = "  -4.04%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "352.27"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -4.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.87"
$ws.Range("E21").Value = "  -1.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.88"
$ws.Range("E22").Value = "  -6.03%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.66"
$ws.Range("E24").Value = "  -2.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.79"
$ws.Range("E25").Value = "  -4.43%  "

$ws.Range("D26").Value = "2.603.13"
$ws.Range("E26").Value = "  -1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.47"
$ws.Range("E27").Value = "  -4.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").Value = "0.0₃0860"
$ws.Range("E29").Value = "  -4.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("E30").Value = "  -6.18%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  -4.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "434.35"
$ws.Range("E32").Value = "  -6.55%  "

$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.70"
$ws.Range("E34").Value = "  -2.44%  "

$ws.Range("B35").Value = "POPCAT"
$ws.Range("C35").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.19"
$ws.Range("E35").Value = "  +115.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.60"
$ws.Range("E36").Value = "  -0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.111"
$ws.Range("E37").Value = "  -4.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.03"
$ws.Range("E40").Value = "  -2.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.310"
$ws.Range("E41").Value = "  -3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.53"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.55"
$ws.Range("E43").Value = "  -3.67%  "

$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.07"
$ws.Range("E44").Value = "  -4.20%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.13"
$ws.Range("E45").Value = "  -4.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.78"
$ws.Range("E46").Value = "  -4.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.500"
$ws.Range("E48").Value = "  -4.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0723"
$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.568"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0920"
$ws.Range("E51").Value = "  -1.61%  "
